# Update countries & provincias Spain
# Applies the COVID-data refresh captured in the commit:
#  - refreshed case numbers for Ucrania (row 28) and Uzbekistan (row 59)
#  - a newly-tracked "Tunez" entry inserted into the country list, which
#    shifts the rows for Guayana Francesa / Maldivas / Tayikistan down by
#    one (with their previous numbers carried along) and drops the old
#    "Tunez" row's stale numbers
#  - "Timor Oriental" / "Santa Lucia" swapped order (numbers identical so
#    only the labels move)
#  - "Islas Malvinas" / "Montserrat" swapped order, carrying their
#    slightly different Casos activos / Muertes numbers with them
#  - refreshed "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Septiembre de 2020 a las 08:25"

# Ucrania (row 28) - refreshed totals
$ws.Range("B28").Value = 175678
$ws.Range("C28").Value = 2966
$ws.Range("D28").Value = 77512
$ws.Range("E28").Value = 94609
$ws.Range("G28").Value = 41
$ws.Range("H28").Value = 3557

# Uzbekistan (row 59) - refreshed totals
$ws.Range("B59").Value = 51235
$ws.Range("C59").Value = 243
$ws.Range("E59").Value = 3535
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 429

# Tunez newly inserted here with fresh numbers; Guayana Francesa,
# Maldivas and Tayikistan each shift down one row keeping their old data
$ws.Range("A99").Value = "Tunez"
$ws.Range("B99").Value = 9736
$ws.Range("C99").Value = 626
$ws.Range("D99").Value = 2386
$ws.Range("E99").Value = 7195
$ws.Range("G99").Value = 17
$ws.Range("H99").Value = 155

$ws.Range("A100").Value = "Guayana Francesa"
$ws.Range("B100").Value = 9692
$ws.Range("D100").Value = 9341
$ws.Range("E100").Value = 286
$ws.Range("H100").Value = 65

$ws.Range("A101").Value = "Maldivas"
$ws.Range("B101").Value = 9649
$ws.Range("D101").Value = 8188
$ws.Range("E101").Value = 1428
$ws.Range("H101").Value = 33

$ws.Range("A102").Value = "Tayikistan"
$ws.Range("B102").Value = 9303
$ws.Range("D102").Value = 8066
$ws.Range("E102").Value = 1164
$ws.Range("H102").Value = 73

# Timor Oriental / Santa Lucia swap order (values were identical)
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

# Islas Malvinas / Montserrat swap order, carrying their own numbers
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
